$d = $word.ActiveDocument

# The edit swaps out the GSM-R radio standard for a second TETRA radio
# ("TETRA + GSM-R" -> "TETRA + TETRA", including in the bullet list, the
# bold table cell, the "components installed" table row, and the
# interfaces table) and lowers the contracted rolling-stock availability
# figure from 99.95% to 99.5% everywhere it is quoted (KPI tables, the
# "why 99.95%?" heading, the compliance/appendix bullet points, etc.).
#
# Both strings occur only inside visible run text (<w:t>) throughout the
# body, each occurrence needing the exact same substitution, so a plain
# document-wide Find & Replace captures every required change.

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("GSM-R", $true, $false, $false, $false, $false, $true, 1, $false, "TETRA", 2)

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("99.95", $true, $false, $false, $false, $false, $true, 1, $false, "99.5", 2)
